$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -7.044999999999999
$ws.Range("C3").Value = -11.358
$ws.Range("D3").Value = -7.216999999999999
$ws.Range("C4").Value = -12.446
$ws.Range("E8").Value = 16.798
$ws.Range("D9").Value = -6.94
$ws.Range("B11").Value = 6.366000000000001
$ws.Range("E11").Value = 16.716
$ws.Range("B12").Value = 4.871
$ws.Range("C14").Value = -13.054
$ws.Range("E14").Value = 16.669
$ws.Range("B15").Value = 4.83
$ws.Range("D15").Value = -8.190000000000001
$ws.Range("E15").Value = 16.392
$ws.Range("E17").Value = 16.771
$ws.Range("D19").Value = -8.099
$ws.Range("D20").Value = -7.825
$ws.Range("D25").Value = -7.95
$ws.Range("C26").Value = -11.443
$ws.Range("E26").Value = 16.71100000000001
$ws.Range("B27").Value = 5.786
$ws.Range("D27").Value = -8.530000000000001
$ws.Range("B28").Value = 6.192
$ws.Range("D28").Value = -7.99
$ws.Range("D30").Value = -7.206
$ws.Range("B31").Value = 5.492
$ws.Range("C31").Value = -12.333
$ws.Range("B32").Value = 6.250999999999999
$ws.Range("D32").Value = -7.593999999999999
$ws.Range("C35").Value = -11.998
$ws.Range("B36").Value = 9.356
$ws.Range("E36").Value = 16.686
$ws.Range("C37").Value = -13.569
$ws.Range("B38").Value = 5.223999999999999
$ws.Range("C39").Value = -12.434
$ws.Range("C40").Value = -12.989
$ws.Range("E42").Value = 16.687
$ws.Range("D44").Value = -7.673999999999999
$ws.Range("C45").Value = -12.776
$ws.Range("B46").Value = 5.696000000000001
$ws.Range("D47").Value = -7.531000000000001
$ws.Range("C52").Value = -11.07
$ws.Range("B54").Value = 5.119
$ws.Range("B55").Value = 4.720999999999999
$ws.Range("B56").Value = 4.614
$ws.Range("C57").Value = -13.646
$ws.Range("D58").Value = -7.974000000000001
$ws.Range("D62").Value = -7.869
$ws.Range("E64").Value = 17.396
$ws.Range("B67").Value = 5.212000000000001
$ws.Range("E68").Value = 17.037
$ws.Range("B69").Value = 5.212
$ws.Range("B72").Value = 5.624
$ws.Range("B73").Value = 8.247
$ws.Range("D77").Value = -7.45
$ws.Range("D78").Value = -8.076999999999998
$ws.Range("E79").Value = 17.394
$ws.Range("C81").Value = -13.035
$ws.Range("B83").Value = 5.007000000000001
$ws.Range("C83").Value = -13.607
$ws.Range("D84").Value = -8.019
$ws.Range("B86").Value = 5.037999999999999
$ws.Range("D89").Value = -7.363
$ws.Range("E89").Value = 17.126
$ws.Range("B91").Value = 5.6
$ws.Range("D91").Value = -6.856999999999999
$ws.Range("D92").Value = -7.192
$ws.Range("B93").Value = 5.77
$ws.Range("D96").Value = -7.446
$ws.Range("B99").Value = 5.811000000000001
$ws.Range("C100").Value = -12.332
$ws.Range("C102").Value = -13.228
$ws.Range("D102").Value = -7.334999999999999
